$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: wrap a body-XML fragment into a minimal WordOpenXML package so
# it can be fed to Range.InsertXML().
# ---------------------------------------------------------------------
function New-WordXmlPackage([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1) Locate the paragraph right after "Test description" (the first of
#    the 4 empty paragraphs in that table cell) and replace it with the
#    full Exercise 1 / 2 / 3 block, ending in the relocated _GoBack
#    bookmark.
# ---------------------------------------------------------------------
$testDescPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Test description`r") {
        $testDescPara = $p
        break
    }
}
$targetPara = $testDescPara.Next()

$exerciseBlock = '<w:p><w:pPr><w:rPr><w:b/><w:i/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>Exercise 1:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Test 1: </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Result 1: </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Test 2: </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Result 2: </w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:i/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:i/></w:rPr><w:t xml:space="preserve">Exercise </w:t></w:r><w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Test 1: </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Result 1: </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Test 2: </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Result 2: </w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:i/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:i/></w:rPr><w:t xml:space="preserve">Exercise </w:t></w:r><w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Test 1: </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Result 1: </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Test 2: </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Result 2: </w:t></w:r></w:p><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$targetPara.Range.InsertXML((New-WordXmlPackage $exerciseBlock))

# ---------------------------------------------------------------------
# 2) Remove the old _GoBack bookmark that used to sit at the end of the
#    "Comments" paragraph ("...squares, circles or more."). Rebuild the
#    paragraph's runs identically, just without the bookmark markers.
# ---------------------------------------------------------------------
$commentsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*With this, I have changed*") {
        $commentsPara = $p
        break
    }
}

$commentsXml = '<w:p w:rsidR="00080257" w:rsidRPr="00080257" w:rsidRDefault="00080257"><w:r w:rsidRPr="00080257"><w:t xml:space="preserve">With this, I have changed my </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00080257"><w:rPr><w:i/></w:rPr><w:t>setLength</w:t></w:r><w:r w:rsidRPr="00080257"><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00080257"><w:t xml:space="preserve">) method to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00080257"><w:rPr><w:i/></w:rPr><w:t>setSize</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>()</w:t></w:r><w:r w:rsidRPr="00080257"><w:t xml:space="preserve"> so that the name is more appropriate for all shape objects, may they be squares, circles or mor</w:t></w:r><w:r><w:t>e.</w:t></w:r></w:p>'

$commentsPara.Range.InsertXML((New-WordXmlPackage $commentsXml))

# ---------------------------------------------------------------------
# 3) Add <w:lastRenderedPageBreak/> right before "Possible improvements".
# ---------------------------------------------------------------------
$possibleImprovementsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Possible improvements`r") {
        $possibleImprovementsPara = $p
        break
    }
}

$possibleImprovementsXml = '<w:p w:rsidR="0056782F" w:rsidRDefault="0056782F"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>Possible improvements</w:t></w:r></w:p>'

$possibleImprovementsPara.Range.InsertXML((New-WordXmlPackage $possibleImprovementsXml))

# ---------------------------------------------------------------------
# 4) Remove <w:lastRenderedPageBreak/> that used to precede "Extra credit".
# ---------------------------------------------------------------------
$extraCreditPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Extra credit`r") {
        $extraCreditPara = $p
        break
    }
}

$extraCreditXml = '<w:p w:rsidR="0056782F" w:rsidRDefault="00381807"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Extra credit</w:t></w:r></w:p>'

$extraCreditPara.Range.InsertXML((New-WordXmlPackage $extraCreditXml))

Write-Host "Done"
